# 3.4 Manage Professor Data — edits per commit "Made edits to Section 3.4"
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Replace the "store faculty information..." paragraph's text with
# the new short summary. This paragraph currently sits right after the
# "3.1.1 Description" heading; we rewrite its text in place (keeping its
# existing sz24/szCs24 run + paragraph-mark formatting) and will relocate the
# heading paragraph to appear after it in a later step.
# ---------------------------------------------------------------------------
$descParaIndex = 4
$pDesc = $d.Paragraphs.Item($descParaIndex)
$descRange = $d.Range($pDesc.Range.Start, $pDesc.Range.End - 1)
$descRange.Text = "Creating and editing professor information by semester."

# Add back the (empty / collapsed) "_GoBack" bookmark immediately after that
# run, inside the same paragraph.
$pDesc2 = $d.Paragraphs.Item($descParaIndex)
$goBackPos = $pDesc2.Range.End - 1
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# ---------------------------------------------------------------------------
# Step 2: Delete the old intro paragraph "Upon creating a new semester,
# importing to an existing semester or editing an existing semester."
# ---------------------------------------------------------------------------
$introPara = $d.Paragraphs.Item(2)
$introPara.Range.Delete()

# ---------------------------------------------------------------------------
# Step 3: Move the "3.1.1 Description" heading paragraph (now index 2) to
# just after the rewritten summary paragraph (now index 2 -> heading is 2,
# summary is 3... after deleting paragraph 2, indices shift up by one, so:
#   1 = Heading2 "3.4 Manage Professor Data"
#   2 = Heading3 "3.1.1 Description"
#   3 = "Creating and editing professor information by semester."
# We cut the heading paragraph and paste it right after paragraph 3.
# ---------------------------------------------------------------------------
$headingPara = $d.Paragraphs.Item(2)
$headingPara.Range.Select()
$word.Selection.Cut()

$summaryPara = $d.Paragraphs.Item(2)
$pastePoint = $d.Range($summaryPara.Range.End, $summaryPara.Range.End)
$pastePoint.Select()
$word.Selection.Paste()

# Re-create the "_Toc479328536" bookmark around the heading text, which the
# cut/paste did not preserve.
$headingParaNew = $d.Paragraphs.Item(3)
$hStart = $headingParaNew.Range.Start
$hEnd = $headingParaNew.Range.End - 1
$hRange = $d.Range($hStart, $hEnd)
$d.Bookmarks.Add("_Toc479328536", $hRange)

# ---------------------------------------------------------------------------
# Step 4: Insert the new, longer description paragraph right after the
# heading (it replaces the paragraph that used to hold the "store faculty
# information..." text but is now a brand-new paragraph because that
# paragraph's text was already consumed/rewritten in Step 1).
# ---------------------------------------------------------------------------
$headingParaFinal = $d.Paragraphs.Item(3)
$headingParaFinal.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(4)
$newPara.Style = $d.Styles.Item("Normal")
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newRange.Text = "An Office Administrator who has been authenticated can view and edit professor information by semester.  An Office Administrator can also add a new professor to the semester.   "
$newRange.Font.Size = 12
$newRange.Font.SizeBi = 12

# ---------------------------------------------------------------------------
# Step 5: Remove the old "_GoBack" bookmark that used to sit between "i" and
# "nformation." near the end of the document.
# ---------------------------------------------------------------------------
Write-Output "Done"
